$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new rows at 507/508 for the new OBInternalStatementFeeType1Code values
$ws.Rows.Item(507).Insert()
$ws.Rows.Item(508).Insert()

$ws.Rows.Item(507).RowHeight = 17
$ws.Rows.Item(508).RowHeight = 17

$ws.Cells.Item(507,1).Value = "OBInternalStatementFeeType1Code"
$ws.Cells.Item(507,2).Value = "UK.OBIE.InstalmentPlan"
$ws.Cells.Item(507,3).Value = "UK.OBIE.InstalmentPlan"
$ws.Cells.Item(507,4).Value = "Instalment plan fees charged during the statement period."

$ws.Cells.Item(508,1).Value = "OBInternalStatementFeeType1Code"
$ws.Cells.Item(508,2).Value = "UK.OBIE.ReturnedPayment"
$ws.Cells.Item(508,3).Value = "UK.OBIE.ReturnedPayment"
$ws.Cells.Item(508,4).Value = "Returned payment fees charged during the statement period."

# Resize the table to include the new rows
$tbl = $ws.ListObjects.Item(1)
$newRange = $ws.Range("A1:F714")
$tbl.Resize($newRange)

# Restore row heights that Excel recalculated on the two unrelated rows
$ws.Rows.Item(349).RowHeight = 68
$ws.Rows.Item(379).RowHeight = 34
$ws.Rows.Item(713).RowHeight = 34
$ws.Rows.Item(714).RowHeight = 34

# Update sheet view to match target state
$ws.Application.ActiveWindow.ScrollRow = 481
$ws.Range("A486").Select()

Write-Host "Edit applied"
